$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextCell 'D2' '246.32'
Set-TextCell 'D3' '22.75'
Set-TextCell 'D5' '0.05731'
Set-TextCell 'D7' '0.8097'
Set-TextCell 'D8' '0.8858'
Set-TextCell 'B9' 'WazirX'
Set-TextCell 'C9' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell 'D9' '0.1426'
Set-TextCell 'E9' '8WazirXWRX'
Set-TextCell 'B10' 'MandalaExchangeToken'
Set-TextCell 'C10' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell 'D10' '0.07369'
Set-TextCell 'E10' '9MandalaExchangeTokenMDX'
Set-TextCell 'B11' 'LiechtensteinCryptoassetsExchange'
Set-TextCell 'C11' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell 'D11' '0.03035'
Set-TextCell 'E11' '10LiechtensteinCryptoassetsExchangeLCX'
Set-TextCell 'B12' 'BitrueCoin'
Set-TextCell 'C12' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell 'D12' '0.03114'
Set-TextCell 'E12' '11BitrueCoinBTR'
Set-TextCell 'B13' 'BitMartToken'
Set-TextCell 'C13' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell 'D13' '0.09395'
Set-TextCell 'E13' '12BitMartTokenBMX'
Set-TextCell 'B14' 'MCDex'
Set-TextCell 'C14' 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextCell 'D14' '3.932'
Set-TextCell 'E14' '13MCDexMCB'
Set-TextCell 'B15' 'BitForexToken'
Set-TextCell 'C15' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell 'D15' '0.001582'
Set-TextCell 'E15' '14BitForexTokenBF'
Set-TextCell 'B16' 'CoinExToken'
Set-TextCell 'C16' 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextCell 'D16' '0.04818'
Set-TextCell 'E16' '15CoinExTokenCET'
Set-TextCell 'B17' 'One'
Set-TextCell 'C17' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextCell 'D17' '0.0005852'
Set-TextCell 'E17' '16OneONE'
Set-TextCell 'D19' '0.005104'
Set-TextCell 'D20' '0.0009968'
Set-TextCell 'D22' '3.747'
Set-TextCell 'D23' '6.310'
Set-TextCell 'D24' '2.189'
Set-TextCell 'B27' 'UpBots'
Set-TextCell 'C27' 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
Set-TextCell 'D27' '0.0002992'
Set-TextCell 'E27' '26UpBotsUBXTWorstin24h'
Set-TextCell 'B28' 'Spectre.aiUtilityToken'
Set-TextCell 'C28' 'https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut'
Set-TextCell 'E28' '27Spectre.aiUtilityTokenSXUT'
Set-TextCell 'B29' 'LegolasExchange'
Set-TextCell 'C29' 'https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo'
Set-TextCell 'E29' '28LegolasExchangeLGO'
Set-TextCell 'B30' 'BitZToken'
Set-TextCell 'C30' 'https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz'
Set-TextCell 'E30' '29BitZTokenBZ'
Set-TextCell 'B31' 'Birake'
Set-TextCell 'C31' 'https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir'
Set-TextCell 'E31' '30BirakeBIR'
Set-TextCell 'B32' 'ZBToken'
Set-TextCell 'C32' 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextCell 'E32' '31ZBTokenZB'
Set-TextCell 'B33' 'NashExchange'
Set-TextCell 'C33' 'https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex'
Set-TextCell 'E33' '32NashExchangeNEX'
Set-TextCell 'B34' 'AAXToken'
Set-TextCell 'C34' 'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab'
Set-TextCell 'E34' '33AAXTokenAAB'
Set-TextCell 'B35' 'CenX'
Set-TextCell 'C35' 'https://coinranking.com/coin/V4XJUvLQb+cenx-cenx'
Set-TextCell 'E35' '34CenXCENX'
Set-TextCell 'B36' 'BNIXToken'
Set-TextCell 'C36' 'https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix'
Set-TextCell 'E36' '35BNIXTokenBNIX'
Set-TextCell 'D40' '0.03909'
Set-TextCell 'D41' '0.006727'
Set-TextCell 'D42' '0.1069'
Set-TextCell 'D43' '0.002530'
Set-TextCell 'D44' '0.007492'
Set-TextCell 'D45' '0.00005636'
Set-TextCell 'D47' '0.6003'
Set-TextCell 'E47' '46CoinbaseStockTokenCOINBestin24h'
Set-TextCell 'D48' '0.1748'
Set-TextCell 'E48' '47BOLOBOLO'
Set-TextCell 'D49' '0.00002101'

Write-Host "Applied 88 cell updates"
